$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Fecha" header in K1, matching the style of the other header cells
$ws.Range("K1").Value = "Fecha"
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill K2:K41 with the date string "06 08 24" for each data row
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 11).Value = "06 08 24"
}
